$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 holds the "publish date" labels for each quarterly column.
# The three most recent dated columns were re-published under a new date.
$ws.Range("I9").Value = "1402-01-28 (5)"
$ws.Range("J9").Value = "1402-01-28 (8)"
$ws.Range("M9").Value = "1402-01-28 (3)"

# Column J (quarter "1401-10-27" -> re-dated "1402-01-28 (8)") values were
# recalculated with the updated read_price algorithm.
$ws.Range("J11").Value = 25152658
$ws.Range("J12").Value = -20998309
$ws.Range("J13").Value = 4154349
$ws.Range("J17").Value = 3537092
$ws.Range("J20").Value = 3475276
$ws.Range("J22").Value = 2992161
$ws.Range("J24").Value = 2992161
$ws.Range("J27").Value = 1496

# Row 25 ("سود هر سهم پس از کسر مالیات") previously had placeholder "-" text
# in columns I, J and M; these are now populated with computed EPS values.
$ws.Range("I25").Value = 2085
$ws.Range("J25").Value = 1496
$ws.Range("M25").Value = 2951
